$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.597.53'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '3.729.67'
$ws.Range("E3").Value = '  -0.05%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = "'613.18"
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("D6").Value = "'178.34"
$ws.Range("E6").Value = '  +0.69%  '
$ws.Range("D7").Value = '3.723.59'
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = '  -2.45%  '
$ws.Range("E10").Value = '  -1.37%  '
$ws.Range("D11").Value = "'6.54"
$ws.Range("E11").Value = '  +2.34%  '
$ws.Range("E12").Value = '  -4.02%  '
$ws.Range("D13").Value = "'39.78"
$ws.Range("E13").Value = '  -2.29%  '
$ws.Range("E14").Value = '  -1.28%  '
$ws.Range("D15").Value = '4.357.76'
$ws.Range("E15").Value = '  +0.12%  '
$ws.Range("D16").Value = '3.723.98'
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("D17").Value = '69.699.33'
$ws.Range("E17").Value = '  -0.33%  '
$ws.Range("E18").Value = '  -2.58%  '
$ws.Range("D19").Value = "'7.46"
$ws.Range("E19").Value = '  -1.73%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = "'16.33"
$ws.Range("E20").Value = '  -2.16%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = "'500.86"
$ws.Range("E21").Value = '  -2.93%  '
$ws.Range("D22").Value = "'9.12"
$ws.Range("E22").Value = '  -2.73%  '
$ws.Range("D23").Value = "'0.719"
$ws.Range("E23").Value = '  -1.16%  '
$ws.Range("E24").Value = '  +4.91%  '
$ws.Range("D25").Value = "'85.99"
$ws.Range("E25").Value = '  -2.64%  '
$ws.Range("D26").Value = "'11.46"
$ws.Range("E26").Value = '  +4.65%  '
$ws.Range("D27").Value = "'12.86"
$ws.Range("E27").Value = '  -5.25%  '
$ws.Range("D28").Value = "'0.0000135"
$ws.Range("E28").Value = '  +6.80%  '
$ws.Range("E29").Value = '  +0.29%  '
$ws.Range("E30").Value = '  -1.98%  '
$ws.Range("E31").Value = '  +2.40%  '
$ws.Range("D32").Value = "'7.99"
$ws.Range("E32").Value = '  +1.72%  '
$ws.Range("D33").Value = "'30.24"
$ws.Range("E33").Value = '  -3.67%  '
$ws.Range("E34").Value = '  -2.09%  '
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = '  -0.18%  '
$ws.Range("E36").Value = '  +0.91%  '
$ws.Range("D37").Value = "'6.08"
$ws.Range("E37").Value = '  -2.13%  '
$ws.Range("D38").Value = "'0.350"
$ws.Range("E38").Value = '  +2.63%  '
$ws.Range("E39").Value = '  +3.38%  '
$ws.Range("D40").Value = "'3.06"
$ws.Range("E40").Value = '  +12.76%  '
$ws.Range("D41").Value = "'2.06"
$ws.Range("E41").Value = '  -5.63%  '
$ws.Range("D42").Value = "'45.81"
$ws.Range("E42").Value = '  +3.50%  '
$ws.Range("D43").Value = "'439.28"
$ws.Range("E43").Value = '  +4.26%  '
$ws.Range("E44").Value = '  -3.11%  '
$ws.Range("D45").Value = "'8.52"
$ws.Range("E45").Value = '  -3.47%  '
$ws.Range("D46").Value = '2.950.59'
$ws.Range("E46").Value = '  -4.05%  '
$ws.Range("E47").Value = '  -1.28%  '
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = '  -0.04%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = "'138.31"
$ws.Range("E49").Value = '  +1.82%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = "'27.08"
$ws.Range("E50").Value = '  -2.90%  '
$ws.Range("D51").Value = "'2.48"
$ws.Range("E51").Value = '  -1.72%  '

# Reset style on quote-prefixed numeric-text cells to avoid residual quotePrefix/style formatting
foreach ($addr in @("D5","D6","D9","D11","D13","D19","D20","D21","D22","D23","D25","D26","D27","D28","D32","D33","D35","D37","D38","D40","D41","D42","D43","D45","D48","D49","D50","D51")) {
  $ws.Range($addr).Style = "Normal"
}
